$d = $word.ActiveDocument
$d.Content.Find.Execute("67+10=", $true, $false, $false, $false, $false, $true, 1, $false, "19+36=", 2) | Out-Null
$d.Content.Find.Execute("73-67=", $true, $false, $false, $false, $false, $true, 1, $false, "14+9=", 2) | Out-Null
$d.Content.Find.Execute("72-49=", $true, $false, $false, $false, $false, $true, 1, $false, "28+6=", 2) | Out-Null
$d.Content.Find.Execute("28+23=", $true, $false, $false, $false, $false, $true, 1, $false, "36+37=", 2) | Out-Null
$d.Content.Find.Execute("8+62=", $true, $false, $false, $false, $false, $true, 1, $false, "47-28=", 2) | Out-Null
$d.Content.Find.Execute("92-0=", $true, $false, $false, $false, $false, $true, 1, $false, "51+2=", 2) | Out-Null
$d.Content.Find.Execute("95-75=", $true, $false, $false, $false, $false, $true, 1, $false, "38+9=", 2) | Out-Null
$d.Content.Find.Execute("91-0=", $true, $false, $false, $false, $false, $true, 1, $false, "24+34=", 2) | Out-Null
$d.Content.Find.Execute("69-32=", $true, $false, $false, $false, $false, $true, 1, $false, "71-48=", 2) | Out-Null
$d.Content.Find.Execute("11+10=", $true, $false, $false, $false, $false, $true, 1, $false, "5+47=", 2) | Out-Null
$d.Content.Find.Execute("44-14=", $true, $false, $false, $false, $false, $true, 1, $false, "75+6=", 2) | Out-Null
$d.Content.Find.Execute("89-8=", $true, $false, $false, $false, $false, $true, 1, $false, "51+40=", 2) | Out-Null
$d.Content.Find.Execute("35+0=", $true, $false, $false, $false, $false, $true, 1, $false, "79-20=", 2) | Out-Null
$d.Content.Find.Execute("75-6=", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=", 2) | Out-Null
$d.Content.Find.Execute("12+80=", $true, $false, $false, $false, $false, $true, 1, $false, "40+18=", 2) | Out-Null
$d.Content.Find.Execute("88-83=", $true, $false, $false, $false, $false, $true, 1, $false, "76-49=", 2) | Out-Null
$d.Content.Find.Execute("36+3=", $true, $false, $false, $false, $false, $true, 1, $false, "38-29=", 2) | Out-Null
$d.Content.Find.Execute("74-10=", $true, $false, $false, $false, $false, $true, 1, $false, "93-22=", 2) | Out-Null
$d.Content.Find.Execute("21+10=", $true, $false, $false, $false, $false, $true, 1, $false, "85-42=", 2) | Out-Null
$d.Content.Find.Execute("84-77=", $true, $false, $false, $false, $false, $true, 1, $false, "54+33=", 2) | Out-Null
$d.Content.Find.Execute("51+48=", $true, $false, $false, $false, $false, $true, 1, $false, "1+57=", 2) | Out-Null
$d.Content.Find.Execute("60+34=", $true, $false, $false, $false, $false, $true, 1, $false, "91-75=", 2) | Out-Null
$d.Content.Find.Execute("52+21=", $true, $false, $false, $false, $false, $true, 1, $false, "34+58=", 2) | Out-Null
$d.Content.Find.Execute("2+85=", $true, $false, $false, $false, $false, $true, 1, $false, "68-8=", 2) | Out-Null
$d.Content.Find.Execute("46+12=", $true, $false, $false, $false, $false, $true, 1, $false, "10+60=", 2) | Out-Null
$d.Content.Find.Execute("32+54=", $true, $false, $false, $false, $false, $true, 1, $false, "38+32=", 2) | Out-Null
$d.Content.Find.Execute("31-22=", $true, $false, $false, $false, $false, $true, 1, $false, "11+75=", 2) | Out-Null
$d.Content.Find.Execute("42-29=", $true, $false, $false, $false, $false, $true, 1, $false, "70+28=", 2) | Out-Null
$d.Content.Find.Execute("36+5=", $true, $false, $false, $false, $false, $true, 1, $false, "74+13=", 2) | Out-Null
$d.Content.Find.Execute("40-17=", $true, $false, $false, $false, $false, $true, 1, $false, "46+27=", 2) | Out-Null
$d.Content.Find.Execute("29+8=", $true, $false, $false, $false, $false, $true, 1, $false, "98-35=", 2) | Out-Null
$d.Content.Find.Execute("37-11=", $true, $false, $false, $false, $false, $true, 1, $false, "22+45=", 2) | Out-Null
$d.Content.Find.Execute("19+65=", $true, $false, $false, $false, $false, $true, 1, $false, "50-11=", 2) | Out-Null
$d.Content.Find.Execute("29-28=", $true, $false, $false, $false, $false, $true, 1, $false, "70-39=", 2) | Out-Null
$d.Content.Find.Execute("92-4=", $true, $false, $false, $false, $false, $true, 1, $false, "6+40=", 2) | Out-Null
$d.Content.Find.Execute("15+8=", $true, $false, $false, $false, $false, $true, 1, $false, "83-35=", 2) | Out-Null
$d.Content.Find.Execute("73-50=", $true, $false, $false, $false, $false, $true, 1, $false, "24-5=", 2) | Out-Null
$d.Content.Find.Execute("16+50=", $true, $false, $false, $false, $false, $true, 1, $false, "62+9=", 2) | Out-Null
$d.Content.Find.Execute("44-0=", $true, $false, $false, $false, $false, $true, 1, $false, "59-11=", 2) | Out-Null
$d.Content.Find.Execute("16+21=", $true, $false, $false, $false, $false, $true, 1, $false, "2+89=", 2) | Out-Null
$d.Content.Find.Execute("15+60=", $true, $false, $false, $false, $false, $true, 1, $false, "13+29=", 2) | Out-Null
$d.Content.Find.Execute("26+73=", $true, $false, $false, $false, $false, $true, 1, $false, "83-30=", 2) | Out-Null
$d.Content.Find.Execute("21+56=", $true, $false, $false, $false, $false, $true, 1, $false, "95-86=", 2) | Out-Null
$d.Content.Find.Execute("80-66=", $true, $false, $false, $false, $false, $true, 1, $false, "80-6=", 2) | Out-Null
$d.Content.Find.Execute("78-34=", $true, $false, $false, $false, $false, $true, 1, $false, "75-0=", 2) | Out-Null
$d.Content.Find.Execute("51+26=", $true, $false, $false, $false, $false, $true, 1, $false, "12+35=", 2) | Out-Null
$d.Content.Find.Execute("42-39=", $true, $false, $false, $false, $false, $true, 1, $false, "5+48=", 2) | Out-Null
$d.Content.Find.Execute("6-3=", $true, $false, $false, $false, $false, $true, 1, $false, "46+43=", 2) | Out-Null
$d.Content.Find.Execute("13+9=", $true, $false, $false, $false, $false, $true, 1, $false, "60-27=", 2) | Out-Null
$d.Content.Find.Execute("78-51=", $true, $false, $false, $false, $false, $true, 1, $false, "19+76=", 2) | Out-Null
$d.Content.Find.Execute("6+78=", $true, $false, $false, $false, $false, $true, 1, $false, "48+20=", 2) | Out-Null
$d.Content.Find.Execute("80+3=", $true, $false, $false, $false, $false, $true, 1, $false, "82-72=", 2) | Out-Null
$d.Content.Find.Execute("86-30=", $true, $false, $false, $false, $false, $true, 1, $false, "69-23=", 2) | Out-Null
$d.Content.Find.Execute("96-34=", $true, $false, $false, $false, $false, $true, 1, $false, "98-86=", 2) | Out-Null
$d.Content.Find.Execute("24-19=", $true, $false, $false, $false, $false, $true, 1, $false, "70+25=", 2) | Out-Null
$d.Content.Find.Execute("34+25=", $true, $false, $false, $false, $false, $true, 1, $false, "43-22=", 2) | Out-Null
$d.Content.Find.Execute("61-26=", $true, $false, $false, $false, $false, $true, 1, $false, "67-37=", 2) | Out-Null
$d.Content.Find.Execute("22-3=", $true, $false, $false, $false, $false, $true, 1, $false, "51+18=", 2) | Out-Null
$d.Content.Find.Execute("19+48=", $true, $false, $false, $false, $false, $true, 1, $false, "48-16=", 2) | Out-Null
$d.Content.Find.Execute("7+30=", $true, $false, $false, $false, $false, $true, 1, $false, "74-38=", 2) | Out-Null
$d.Content.Find.Execute("31+5=", $true, $false, $false, $false, $false, $true, 1, $false, "83-43=", 2) | Out-Null
$d.Content.Find.Execute("46-25=", $true, $false, $false, $false, $false, $true, 1, $false, "57+14=", 2) | Out-Null
$d.Content.Find.Execute("93-20=", $true, $false, $false, $false, $false, $true, 1, $false, "41+10=", 2) | Out-Null
$d.Content.Find.Execute("77-22=", $true, $false, $false, $false, $false, $true, 1, $false, "44+4=", 2) | Out-Null
$d.Content.Find.Execute("44+42=", $true, $false, $false, $false, $false, $true, 1, $false, "66+21=", 2) | Out-Null
$d.Content.Find.Execute("92+4=", $true, $false, $false, $false, $false, $true, 1, $false, "78-55=", 2) | Out-Null
$d.Content.Find.Execute("96-62=", $true, $false, $false, $false, $false, $true, 1, $false, "20-15=", 2) | Out-Null
$d.Content.Find.Execute("23+28=", $true, $false, $false, $false, $false, $true, 1, $false, "66-48=", 2) | Out-Null
$d.Content.Find.Execute("19+80=", $true, $false, $false, $false, $false, $true, 1, $false, "60-4=", 2) | Out-Null
$d.Content.Find.Execute("10+47=", $true, $false, $false, $false, $false, $true, 1, $false, "80-50=", 2) | Out-Null
$d.Content.Find.Execute("25-8=", $true, $false, $false, $false, $false, $true, 1, $false, "87-44=", 2) | Out-Null
$d.Content.Find.Execute("57+24=", $true, $false, $false, $false, $false, $true, 1, $false, "72+3=", 2) | Out-Null
$d.Content.Find.Execute("83-9=", $true, $false, $false, $false, $false, $true, 1, $false, "21+71=", 2) | Out-Null
$d.Content.Find.Execute("80-18=", $true, $false, $false, $false, $false, $true, 1, $false, "71-43=", 2) | Out-Null
$d.Content.Find.Execute("88-62=", $true, $false, $false, $false, $false, $true, 1, $false, "23+30=", 2) | Out-Null
$d.Content.Find.Execute("74+24=", $true, $false, $false, $false, $false, $true, 1, $false, "71-0=", 2) | Out-Null
$d.Content.Find.Execute("85+3=", $true, $false, $false, $false, $false, $true, 1, $false, "82-58=", 2) | Out-Null
$d.Content.Find.Execute("56-17=", $true, $false, $false, $false, $false, $true, 1, $false, "98-17=", 2) | Out-Null
$d.Content.Find.Execute("70-35=", $true, $false, $false, $false, $false, $true, 1, $false, "10+83=", 2) | Out-Null
$d.Content.Find.Execute("37+23=", $true, $false, $false, $false, $false, $true, 1, $false, "15+80=", 2) | Out-Null
$d.Content.Find.Execute("73-60=", $true, $false, $false, $false, $false, $true, 1, $false, "17+4=", 2) | Out-Null
$d.Content.Find.Execute("90-76=", $true, $false, $false, $false, $false, $true, 1, $false, "57-50=", 2) | Out-Null
$d.Content.Find.Execute("69-13=", $true, $false, $false, $false, $false, $true, 1, $false, "36-10=", 2) | Out-Null
$d.Content.Find.Execute("96-7=", $true, $false, $false, $false, $false, $true, 1, $false, "13+55=", 2) | Out-Null
$d.Content.Find.Execute("88-67=", $true, $false, $false, $false, $false, $true, 1, $false, "17+29=", 2) | Out-Null
$d.Content.Find.Execute("79-73=", $true, $false, $false, $false, $false, $true, 1, $false, "79+17=", 2) | Out-Null
$d.Content.Find.Execute("29+13=", $true, $false, $false, $false, $false, $true, 1, $false, "32+3=", 2) | Out-Null
$d.Content.Find.Execute("99-86=", $true, $false, $false, $false, $false, $true, 1, $false, "12+4=", 2) | Out-Null
$d.Content.Find.Execute("77-70=", $true, $false, $false, $false, $false, $true, 1, $false, "82+6=", 2) | Out-Null
$d.Content.Find.Execute("17-0=", $true, $false, $false, $false, $false, $true, 1, $false, "40+31=", 2) | Out-Null
$d.Content.Find.Execute("49-35=", $true, $false, $false, $false, $false, $true, 1, $false, "11+77=", 2) | Out-Null
$d.Content.Find.Execute("99-36=", $true, $false, $false, $false, $false, $true, 1, $false, "33-31=", 2) | Out-Null
$d.Content.Find.Execute("42-27=", $true, $false, $false, $false, $false, $true, 1, $false, "78-15=", 2) | Out-Null
$d.Content.Find.Execute("11+53=", $true, $false, $false, $false, $false, $true, 1, $false, "11+74=", 2) | Out-Null
$d.Content.Find.Execute("54-40=", $true, $false, $false, $false, $false, $true, 1, $false, "0+47=", 2) | Out-Null
$d.Content.Find.Execute("37-16=", $true, $false, $false, $false, $false, $true, 1, $false, "66+32=", 2) | Out-Null
$d.Content.Find.Execute("54-33=", $true, $false, $false, $false, $false, $true, 1, $false, "99-89=", 2) | Out-Null
$d.Content.Find.Execute("35+7=", $true, $false, $false, $false, $false, $true, 1, $false, "77-72=", 2) | Out-Null
$d.Content.Find.Execute("24+1=", $true, $false, $false, $false, $false, $true, 1, $false, "77-32=", 2) | Out-Null
$d.Content.Find.Execute("84+12=", $true, $false, $false, $false, $false, $true, 1, $false, "21+49=", 2) | Out-Null
